$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    if ($null -eq $value) { return }
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$changes = @(
    @{ Row = 2; D = '320.71'; E = '-3.21%'; G = '11' },
    @{ Row = 3; D = '42.45'; E = '-6.36%'; G = '11' },
    @{ Row = 4; D = '5.190'; E = '-6.68%'; G = '11' },
    @{ Row = 5; D = '0.08171'; E = '-2.16%'; G = '11' },
    @{ Row = 6; D = '4.310'; E = '-3.23%'; G = '11' },
    @{ Row = 7; D = '1.812'; E = '-13.24%'; G = '11' },
    @{ Row = 8; D = '0.9340'; E = '-5.00%'; G = '11' },
    @{ Row = 9; D = '0.1109'; E = '-7.69%'; G = '11' },
    @{ Row = 10; D = '0.1869'; E = '-2.38%'; G = '11' },
    @{ Row = 11; D = '0.09425'; E = '-4.26%'; G = '11' },
    @{ Row = 12; D = '0.04692'; E = '0.58%'; G = '11' },
    @{ Row = 13; D = '7.407'; E = '-28.17%'; G = '11' },
    @{ Row = 14; D = '0.1059'; E = '0.11%'; G = '11' },
    @{ Row = 15; D = '0.001299'; E = '0.77%'; G = '11' },
    @{ Row = 16; D = '0.005705'; E = '-4.12%'; G = '11' },
    @{ Row = 17; D = '3.356'; E = '-1.12%'; G = '11' },
    @{ Row = 18; D = $null; E = '-0.20%'; G = '11' },
    @{ Row = 19; D = '0.3376'; E = '0.80%'; G = '11' },
    @{ Row = 20; D = '0.1389'; E = '2.58%'; G = '11' },
    @{ Row = 21; D = '0.2547'; E = '-0.70%'; G = '11' },
    @{ Row = 22; D = '0.04156'; E = $null; G = '11' },
    @{ Row = 23; D = '0.001246'; E = '-3.70%'; G = '11' },
    @{ Row = 24; D = '0.004325'; E = '-5.43%'; G = '11' },
    @{ Row = 25; D = '0.0001202'; E = '-7.69%'; G = '11' },
    @{ Row = 26; D = '0.0002983'; E = '-20.35%'; G = '11' },
    @{ Row = 27; D = $null; E = $null; G = '11' },
    @{ Row = 28; D = $null; E = $null; G = '11' },
    @{ Row = 29; D = $null; E = $null; G = '11' },
    @{ Row = 30; D = $null; E = $null; G = '11' },
    @{ Row = 31; D = $null; E = $null; G = '11' },
    @{ Row = 32; D = $null; E = $null; G = '11' },
    @{ Row = 33; D = $null; E = $null; G = '11' },
    @{ Row = 34; D = $null; E = $null; G = '11' },
    @{ Row = 35; D = $null; E = $null; G = '11' },
    @{ Row = 36; D = $null; E = $null; G = '11' },
    @{ Row = 37; D = $null; E = $null; G = '11' },
    @{ Row = 38; D = '0.02704'; E = '0.65%'; G = '11' },
    @{ Row = 39; D = $null; E = '-3.27%'; G = '11' },
    @{ Row = 40; D = '0.008054'; E = '2.09%'; G = '11' },
    @{ Row = 41; D = '0.1397'; E = '-2.32%'; G = '11' },
    @{ Row = 42; D = '0.006560'; E = '-12.69%'; G = '11' },
    @{ Row = 43; D = '0.002089'; E = '-1.48%'; G = '11' },
    @{ Row = 44; D = '0.008264'; E = '-8.30%'; G = '11' },
    @{ Row = 45; D = '0.3482'; E = '2.40%'; G = '11' },
    @{ Row = 46; D = '0.00006916'; E = '-2.34%'; G = '11' },
    @{ Row = 47; D = '0.00000000751'; E = '0.00%'; G = '11' },
    @{ Row = 48; D = '0.003370'; E = '-4.49%'; G = '11' },
    @{ Row = 49; D = '0.003535'; E = '0.00%'; G = '11' },
    @{ Row = 50; D = '0.00002103'; E = '0.00%'; G = '11' },
    @{ Row = 51; D = '0.0002003'; E = '0.00%'; G = '11' }
)

foreach ($item in $changes) {
    $r = $item.Row
    if ($null -ne $item.D) { Set-TextValue $ws.Cells.Item($r, 4) $item.D }
    if ($null -ne $item.E) { Set-TextValue $ws.Cells.Item($r, 5) $item.E }
    if ($null -ne $item.G) { Set-TextValue $ws.Cells.Item($r, 7) $item.G }
}